$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the status text everywhere it appears: "Ready for handoff"
#    becomes "Handed back: in sync with en-US" (files have now been handed
#    back from translation and are in sync with en-US).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B2").Value = $newStatus
$wsZhCn.Range("B3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B2").Value = $newStatus
$wsDeDe.Range("B3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Fill in the "Latest Target File" (E) / "Latest Handback File" (F)
#    columns for rows 2 and 3 on both the zh-cn and de-de sheets, add
#    matching hyperlinks (mirroring the Source File Name / Latest Handoff
#    File hyperlinks already on columns A and C), and stamp the real
#    "Latest Handback DateTime" (G) now that the handback happened.
# ---------------------------------------------------------------------------
function Set-HandbackRow {
    param(
        $ws,
        [int]$row,
        [string]$handbackDateTime
    )

    $aLink = $null
    $cLink = $null
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address($false, $false)
        if ($addr -eq ("A" + $row)) { $aLink = $h }
        if ($addr -eq ("C" + $row)) { $cLink = $h }
    }

    $eCell = $ws.Range("E" + $row)
    $fCell = $ws.Range("F" + $row)

    $eCell.Value = $aLink.TextToDisplay
    $fCell.Value = $cLink.TextToDisplay

    $ws.Hyperlinks.Add($eCell, $aLink.Address, "", "", $aLink.TextToDisplay) | Out-Null
    $ws.Hyperlinks.Add($fCell, $cLink.Address, "", "", $cLink.TextToDisplay) | Out-Null

    $eCell.Style = $ws.Range("A" + $row).Style
    $fCell.Style = $ws.Range("C" + $row).Style

    $ws.Range("G" + $row).Value = $handbackDateTime
}

Set-HandbackRow $wsZhCn 2 "2016-02-16 15:43:15"
Set-HandbackRow $wsZhCn 3 "2016-02-16 15:43:15"

Set-HandbackRow $wsDeDe 2 "2016-02-16 15:43:45"
Set-HandbackRow $wsDeDe 3 "2016-02-16 15:43:45"
